# Apply the scene-category block-order update:
# shuffle the header labels across A1:F1 and the corresponding
# one-hot "block order" matrix in rows 2-7.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: header labels (A1:F1) ---
$ws.Range("A1").Value = "kitchens_1"
$ws.Range("B1").Value = "bedrooms_1"
$ws.Range("C1").Value = "living_rooms_1"
$ws.Range("D1").Value = "living_rooms_2"
$ws.Range("E1").Value = "kitchens_2"
$ws.Range("F1").Value = "bedrooms_2"

# --- Rows 2-7: one-hot matrix values ---
$values = @(
    @(0, 1, 0, 0, 0, 0),
    @(0, 0, 0, 0, 1, 0),
    @(0, 0, 0, 1, 0, 0),
    @(0, 0, 0, 0, 0, 1),
    @(1, 0, 0, 0, 0, 0),
    @(0, 0, 1, 0, 0, 0)
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $rowValues = $values[$i]
    for ($j = 0; $j -lt $rowValues.Length; $j++) {
        $col = $j + 1
        $ws.Cells.Item($row, $col).Value = $rowValues[$j]
    }
}
